$wb = $excel.ActiveWorkbook

# --- Sheet "Patient": insert a new "visibilityStatus" column before the
#     existing "villageId" column (K), shifting villageId -> L and
#     mergedIntoId -> M. Also refresh the sample row's data. ---
$wsPatient = $wb.Worksheets.Item("Patient")

$wsPatient.Columns("K").Insert()
$wsPatient.Cells.Item(1, 11).Value = "visibilityStatus"
$wsPatient.Cells.Item(2, 11).Value = "current"

$wsPatient.Cells.Item(2, 1).Value = "7e34e48c-6482-4e81-9fea-de3a04724635"
$wsPatient.Cells.Item(2, 2).Value = "UTLE519661"
$wsPatient.Cells.Item(2, 3).Value = "Betty"
$wsPatient.Cells.Item(2, 5).Value = "Thompson"
$wsPatient.Cells.Item(2, 6).Value = "Louis"
$wsPatient.Cells.Item(2, 7).Value = 33574

# --- Sheet "Allergy": append a "visibilityStatus" column, every existing
#     reference row is "current". ---
$wsAllergy = $wb.Worksheets.Item("Allergy")

$wsAllergy.Cells.Item(1, 4).Value = "visibilityStatus"
$wsAllergy.Cells.Item(2, 4).Value = "current"
$wsAllergy.Cells.Item(3, 4).Value = "current"

# --- Sheet "Diagnosis": same "visibilityStatus" column addition. ---
$wsDiagnosis = $wb.Worksheets.Item("Diagnosis")

$wsDiagnosis.Cells.Item(1, 4).Value = "visibilityStatus"
$wsDiagnosis.Cells.Item(2, 4).Value = "current"
$wsDiagnosis.Cells.Item(3, 4).Value = "current"
